$wb = $excel.ActiveWorkbook

# ---- Sheet "v_batch_size": add batch sizes 128, 256 (rows 16-17) and blank row 18 ----
$wsBatch = $wb.Worksheets.Item("v_batch_size")

$wsBatch.Range("B16").Value = 6
$wsBatch.Range("C16").Value = 128
$wsBatch.Range("D16").Value = 0.0679
$wsBatch.Range("E16").Value = 0.0622

$wsBatch.Range("B17").Value = 7
$wsBatch.Range("C17").Value = 256
$wsBatch.Range("D17").Value = 0.0349
$wsBatch.Range("E17").Value = 0.0622

$wsBatch.Range("D16:E18").NumberFormat = "0.0000"

$wsBatch.Range("A2").Select()

# ---- Sheet "v_n_bottleneck": add bottleneck sizes 128, 256, 512 (rows 15-17) ----
$wsBottleneck = $wb.Worksheets.Item("v_n_bottleneck")

$wsBottleneck.Range("B15").Value = 5
$wsBottleneck.Range("C15").Value = 128
$wsBottleneck.Range("D15").Value = 0.0652
$wsBottleneck.Range("E15").Value = 0.063

$wsBottleneck.Range("B16").Value = 6
$wsBottleneck.Range("C16").Value = 256
$wsBottleneck.Range("D16").Value = 0.0652
$wsBottleneck.Range("E16").Value = 0.063

$wsBottleneck.Range("B17").Value = 7
$wsBottleneck.Range("C17").Value = 512
$wsBottleneck.Range("D17").Value = 0.0652
$wsBottleneck.Range("E17").Value = 0.063

$wsBottleneck.Range("D15:E17").NumberFormat = "0.0000"

$wsBottleneck.Range("P13").Select()

# ---- Update chart series ranges on "Model" sheet to reflect the new data ----
$wsModel = $wb.Worksheets.Item("Model")

$chart1 = $wsModel.ChartObjects(1).Chart
$chart1.SeriesCollection(1).XValues = $wsBatch.Range("C10:C18")
$chart1.SeriesCollection(1).Values = $wsBatch.Range("D10:D18")
$chart1.SeriesCollection(2).XValues = $wsBatch.Range("C10:C18")
$chart1.SeriesCollection(2).Values = $wsBatch.Range("E10:E18")

$chart2 = $wsModel.ChartObjects(2).Chart
$chart2.SeriesCollection(1).XValues = $wsBottleneck.Range("C10:C18")
$chart2.SeriesCollection(1).Values = $wsBottleneck.Range("D10:D17")
$chart2.SeriesCollection(2).XValues = $wsBottleneck.Range("C10:C18")
$chart2.SeriesCollection(2).Values = $wsBottleneck.Range("E10:E17")
